$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the fund from "Agri Fund" to "SAAS Fund" for all data rows
$ws.Range("A2").Value = "SAAS Fund"
$ws.Range("A3").Value = "SAAS Fund"
$ws.Range("A4").Value = "SAAS Fund"

# Update the active selection to A4
$ws.Range("A4").Select()
